# Apply crypto price/volume updates scraped on Thu May 23 23:54:24 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.877.93"
$ws.Range("E2").Value = "  -1.89%  "

$ws.Range("D3").Value = "'3.759.63"
$ws.Range("E3").Value = "  +0.41%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'598.22"
$ws.Range("E5").Value = "  -2.73%  "

$ws.Range("D6").Value = "'176.16"
$ws.Range("E6").Value = "  -0.41%  "

$ws.Range("D7").Value = "'3.757.46"
$ws.Range("E7").Value = "  +0.41%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  +0.33%  "

$ws.Range("E10").Value = "  -3.93%  "

$ws.Range("D11").Value = "'6.21"
$ws.Range("E11").Value = "  -5.16%  "

$ws.Range("D13").Value = "'38.53"
$ws.Range("E13").Value = "  -3.51%  "

$ws.Range("D14").Value = "'0.0000246"
$ws.Range("E14").Value = "  -2.89%  "

$ws.Range("D15").Value = "'4.385.13"
$ws.Range("E15").Value = "  +0.39%  "

$ws.Range("D16").Value = "'3.755.74"
$ws.Range("E16").Value = "  +0.30%  "

$ws.Range("D17").Value = "'67.760.67"
$ws.Range("E17").Value = "  -2.15%  "

$ws.Range("E18").Value = "  -4.14%  "

$ws.Range("E19").Value = "  -3.40%  "

$ws.Range("D20").Value = "'16.56"
$ws.Range("E20").Value = "  +1.19%  "

$ws.Range("D21").Value = "'491.18"
$ws.Range("E21").Value = "  -1.62%  "

$ws.Range("E22").Value = "  -2.97%  "

$ws.Range("D23").Value = "'0.742"
$ws.Range("E23").Value = "  +2.49%  "

$ws.Range("D24").Value = "'85.41"
$ws.Range("E24").Value = "  -0.42%  "

$ws.Range("E25").Value = "  +11.63%  "

$ws.Range("E26").Value = "  -6.38%  "

$ws.Range("D27").Value = "'12.29"

$ws.Range("D28").Value = "'10.23"
$ws.Range("E28").Value = "  -3.99%  "

$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.13%  "

$ws.Range("D30").Value = "'2.96"
$ws.Range("E30").Value = "  +0.57%  "

$ws.Range("E31").Value = "  -2.98%  "

$ws.Range("D32").Value = "'32.20"
$ws.Range("E32").Value = "  +5.80%  "

$ws.Range("D33").Value = "'7.69"
$ws.Range("E33").Value = "  -3.59%  "

$ws.Range("E34").Value = "  -4.04%  "

$ws.Range("D35").Value = "'0.998"
$ws.Range("E35").Value = "  -0.04%  "

$ws.Range("E36").Value = "  -4.47%  "

$ws.Range("E37").Value = "  -5.13%  "

$ws.Range("D38").Value = "'0.134"
$ws.Range("E38").Value = "  -1.94%  "

$ws.Range("E39").Value = "  -5.36%  "

$ws.Range("D40").Value = "'447.89"
$ws.Range("E40").Value = "  +0.45%  "

$ws.Range("D41").Value = "'49.20"
$ws.Range("E41").Value = "  -0.90%  "

$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'2.00"
$ws.Range("E42").Value = "  -2.92%  "

$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "'2.93"
$ws.Range("E43").Value = "  -3.23%  "

$ws.Range("E44").Value = "  -2.37%  "

$ws.Range("D45").Value = "'41.27"
$ws.Range("E45").Value = "  -7.67%  "

$ws.Range("D46").Value = "'2.827.61"
$ws.Range("E46").Value = "  -3.85%  "

$ws.Range("D48").Value = "'138.55"
$ws.Range("E48").Value = "  +0.00%  "

$ws.Range("D49").Value = "'0.0351"
$ws.Range("E49").Value = "  -2.15%  "

$ws.Range("D50").Value = "'26.10"
$ws.Range("E50").Value = "  -4.17%  "

$ws.Range("D51").Value = "'23.73"
$ws.Range("E51").Value = "  +7.73%  "
